$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 9 with a single cell A9 containing "FILLER",
# mirroring the style used by the other column-header cells in column A
# (e.g. A7). Row 8 is intentionally left empty.
$ws.Range("A9").Value = "FILLER"
$ws.Range("A9").Font.Bold = $true
$ws.Range("A9").Font.Color = $ws.Range("A7").Font.Color

# Update the active selection to reflect where the cursor ended up
# after the edit (next empty row below the new data).
$ws.Range("A10").Select()
